# Automatic update of files.
# The underlying data rows (2-13) got reshuffled: each row's species/record
# block (Id, Taxonsorteringsordning, TaxonId, Artnamn, Vetenskapligt namn,
# Auktor, Ost, Nord and the public-comment field) moved to a different row,
# while the location/date/observer columns stayed put (they are identical
# for every row anyway). Row 3 keeps its original content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = 111815515
$ws.Range("Q2").Value = 458161.9437607233
$ws.Range("R2").Value = 7054459.400503729

# --- Row 4 --- (was Tretåig hackspett -> becomes Garnlav) ----------------
$ws.Range("A4").Value = 111815517
$ws.Range("B4").Value = 77515
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("Q4").Value = 458250.8216980004
$ws.Range("R4").Value = 7054375.482693202
$ws.Range("AC4").ClearContents()

# --- Row 5 --- (was Tretåig hackspett -> becomes Granticka) --------------
$ws.Range("A5").Value = 111815514
$ws.Range("B5").Value = 89423
$ws.Range("E5").Value = 5432
$ws.Range("F5").Value = "Granticka"
$ws.Range("G5").Value = "Porodaedalea chrysoloma"
$ws.Range("H5").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("K5").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("Q5").Value = 458153.7808649908
$ws.Range("R5").Value = 7054482.19637617
$ws.Range("AC5").ClearContents()

# --- Row 6 --- (was Granticka -> becomes Tretåig hackspett) --------------
$ws.Range("A6").Value = 111815513
$ws.Range("B6").Value = 56398
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = ""
$ws.Range("Q6").Value = 458173.7327805056
$ws.Range("R6").Value = 7054711.474791372
$ws.Range("AC6").Value = "ringhack gamla"

# --- Row 7 --- (was Garnlav -> becomes Granticka) -------------------------
$ws.Range("A7").Value = 111815516
$ws.Range("B7").Value = 89423
$ws.Range("E7").Value = 5432
$ws.Range("F7").Value = "Granticka"
$ws.Range("G7").Value = "Porodaedalea chrysoloma"
$ws.Range("H7").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q7").Value = 458289.5512131723
$ws.Range("R7").Value = 7054475.069158822

# --- Row 8 ---
$ws.Range("A8").Value = 111815508
$ws.Range("Q8").Value = 458162.4570845839
$ws.Range("R8").Value = 7054329.489790585
$ws.Range("AC8").Value = "ringhack"

# --- Row 9 ---
$ws.Range("A9").Value = 111815507
$ws.Range("Q9").Value = 458151.5539710881
$ws.Range("R9").Value = 7054482.225765129
$ws.Range("AC9").Value = "ringhack gamla"

# --- Row 10 ---
$ws.Range("A10").Value = 111815512
$ws.Range("Q10").Value = 458154.6107204149
$ws.Range("R10").Value = 7054646.336103803
$ws.Range("AC10").Value = "ringhack"

# --- Row 11 ---
$ws.Range("A11").Value = 111815518
$ws.Range("Q11").Value = 458250.901553072
$ws.Range("R11").Value = 7054618.376188213

# --- Row 12 ---
$ws.Range("A12").Value = 111815519
$ws.Range("Q12").Value = 458215.7474518137
$ws.Range("R12").Value = 7054621.063481365

# --- Row 13 --- (was Granticka -> becomes Tretåig hackspett) -------------
$ws.Range("A13").Value = 111815510
$ws.Range("B13").Value = 56398
$ws.Range("E13").Value = 100109
$ws.Range("F13").Value = "Tretåig hackspett"
$ws.Range("G13").Value = "Picoides tridactylus"
$ws.Range("H13").Value = "(Linnaeus, 1758)"
$ws.Range("K13").Value = ""
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("Q13").Value = 458203.7272220219
$ws.Range("R13").Value = 7054385.000644128
$ws.Range("AC13").Value = "ringhack"
